$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new product "Pain relief pills" in column D, row 4
$ws.Range("D4").Value = "Pain relief pills"

# Set column D width to match the diff (14.85546875 chars, closest reachable grid value)
$ws.Columns.Item(4).ColumnWidth = 14

# Update the active selection to E8, matching the post-edit cursor position
$ws.Range("E8").Select()
